# Feature : Sword의 무기 스킬 추가
# Bump the "_range" stat (column Q) from 5 to 7 for the Axe/Sword weapon
# rows (rows 2-9) on the "Entities" sheet, then leave the selection on the
# edited range like the author did before saving.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Entities")

for ($row = 2; $row -le 9; $row++) {
    $ws.Cells.Item($row, 17).Value2 = 7
}

$ws.Range("Q2:Q9").Select()
